# Add biological context fields (organism, anatomy, cell_type) to the
# "Sample" sheet, inserted right after "storage_conditions" (column I) and
# before "parent_sample_id" (old column J).
#
# This shifts the existing parent_sample_id, purity_percentage,
# quality_metrics, id, title, description columns from J:O to M:R, growing
# the sheet's used range from A1:O1 to A1:R1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

# Insert 3 blank columns at J:L, pushing the old J:L (and everything after)
# three columns to the right.
$ws.Range("J1:L1").EntireColumn.Insert()

# Populate the headers for the newly inserted columns.
$ws.Range("J1").Value = "organism"
$ws.Range("K1").Value = "anatomy"
$ws.Range("L1").Value = "cell_type"
